$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws.Range("B9").Value = "Alvearie Team"

# Old "Contact" / "No display for ContactDetail" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The old duplicate "Contact" row (row 11) is removed entirely, shifting remaining rows up
$ws.Rows.Item(11).Delete()
